$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Contribución a los entregables")
$ws1.Range("C2").Value = 1.2
$ws1.Range("C5").Value = 0.8

$ws2 = $wb.Worksheets.Item("Comentarios sobre #DP")
$ws2.Range("A1:C1").Merge()
$ws2.Range("A1").Value = "Manuel Chica Lopez se incorporó mas tarde en el grupo, por esto tiene menos horas en Clockify"
$ws2.Range("A1").WrapText = $true
$ws2.Range("A1").HorizontalAlignment = -4108
$ws2.Rows.Item(1).RowHeight = 101.5
